$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.592.23"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.84"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.37"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4280"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3664"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07271"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8707"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.799.25"
$ws.Range("E12").Value = "  -7.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.424"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.540"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06948"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.37"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008936"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.44"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.538.40"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.176"
$ws.Range("E22").Value = "  +3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.061.99"
$ws.Range("E24").Value = "  -4.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.983"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.44"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.225"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.76"
$ws.Range("E29").Value = "  -5.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.839"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08887"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7614"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.556"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.952"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.143"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05317"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01945"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.813"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1670"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5092"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.617"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.448"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.52"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.10"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06505"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4689"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.621"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.768"
$ws.Range("E51").Value = "  +3.21%  "
